# Rename "strain" input sheets to "strain_log2_expression" (per commit message),
# and update the remembered selection on the renamed "dcin5" sheet.

$wb = $excel.ActiveWorkbook

$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Remember which sheet is currently active so the active tab isn't disturbed -
# selecting a cell on another sheet requires activating it first.
$activeSheetName = $wb.ActiveSheet.Name

$wsDcin5.Activate()
$wsDcin5.Range("G29").Select()

# Restore the original active sheet/tab.
$wb.Worksheets.Item($activeSheetName).Activate()
